$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.908200000000003

$ws.Range("C9").Value = -13.0245
$ws.Range("D9").Value = -8.149600000000001

$ws.Range("C18").Value = -11.68199999999999

$ws.Range("C20").Value = -11.4724

$ws.Range("D23").Value = -8.046999999999997

$ws.Range("D24").Value = -7.480100000000001

$ws.Range("D26").Value = -7.361500000000002

$ws.Range("C27").Value = -12.22849999999999

$ws.Range("D34").Value = -8.170600000000004

$ws.Range("D35").Value = -7.810400000000003

$ws.Range("D48").Value = -7.968699999999999

$ws.Range("D52").Value = -7.808900000000001

$ws.Range("D66").Value = -7.0755

$ws.Range("D67").Value = -7.293399999999997

$ws.Range("C69").Value = -10.6245

$ws.Range("C76").Value = -12.8771

$ws.Range("D80").Value = -8.049399999999999

$ws.Range("C82").Value = -11.928

$ws.Range("D99").Value = -8.231800000000005
